$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for data rows 2-43 is updated from serial date
# 45809 (2025-06-01) to 45810 (2025-06-02).
$ws.Range("C2:C43").Value = 45810
